$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 131 - this shifts the existing rows 131..266
# down to 132..267 (same as the canonical diff, which is a pure row-insert
# followed by population of the newly inserted row with a new weekly
# price entry).
$ws.Rows(131).Insert()

# Populate the newly inserted row 131 with the new data point.
$ws.Range('A131').Value = 4
$ws.Range('B131').Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range('C131').Value = 'Los Lagos'
$ws.Range('D131').Value = 44705
$ws.Range('E131').Value = 10
$ws.Range('F131').Value = 100112003
$ws.Range('G131').Value = 'Ajo'
$ws.Range('H131').Value = 'Chino'
$ws.Range('I131').Value = 'Primera'
$ws.Range('J131').Value = 180
$ws.Range('K131').Value = 21000
$ws.Range('L131').Value = 22000
$ws.Range('M131').Value = 21500
$ws.Range('N131').Value = '$/caja 10 kilos'
$ws.Range('O131').Value = 'China'
$ws.Range('P131').Value = 2150
$ws.Range('Q131').Value = 10
$ws.Range('R131').Value = 'Hortaliza'
